# Update scripts with new TPM values for Vtn-Itga8 LR pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New numeric values (TPM-updated) for rows 2-13 (columns G,H,I,J,M,N,O,P,Q,R,S,T)
$rowData = @{
    2  = @{ G=7.844453333333334;  H=23.53336;           I=0.1489290605659587; J=0.1489290605659588;
            M=1.115279333333333;  N=3.345838;            O=0.08670767158519405; P=0.08670767158519403;
            Q=8.748756683964444;  R=78.73881015568;      S=0.01291329207304462; T=0.01291329207304462 }
    3  = @{ G=7.844453333333334;  H=23.53336;           I=0.1489290605659587; J=0.1489290605659588;
            M=3.484068333333333;  N=10.452205;           O=0.2708697667015328;  P=0.2708697667015328;
            Q=27.33061145097778;  R=245.9755030588;      S=0.04034037989057969; T=0.0403403798905797 }
    4  = @{ G=7.844453333333334;  H=23.53336;           I=0.1489290605659587; J=0.1489290605659588;
            M=8.022733000000001;  N=24.068199;           O=0.6237293899283516;  P=0.6237293899283515;
            Q=62.93395462429334;  R=566.4055916186401;   S=0.09289143208940796; T=0.09289143208940796 }
    5  = @{ G=7.844453333333334;  H=23.53336;           I=0.1489290605659587; J=0.1489290605659588;
            M=0.2404413333333334; N=0.7213240000000001;  O=0.0186931717849216;  P=0.0186931717849216;
            Q=1.886130818737778;  R=16.97517736864;      S=0.00278395651292646; T=0.00278395651292646 }
    6  = @{ G=20.35396833333334;  H=61.06190500000001;  I=0.3864255740794268; J=0.3864255740794268;
            M=1.115279333333333;  N=3.345838;            O=0.08670767158519405; P=0.08670767158519403;
            Q=22.70036023348778;  R=204.30324210139;     S=0.03350606176939901; T=0.03350606176939901 }
    7  = @{ G=20.35396833333334;  H=61.06190500000001;  I=0.3864255740794268; J=0.3864255740794268;
            M=3.484068333333333;  N=10.452205;           O=0.2708697667015328;  P=0.2708697667015328;
            Q=70.91461652783612;  R=638.231548750525;    S=0.1046710050984002;  T=0.1046710050984002 }
    8  = @{ G=20.35396833333334;  H=61.06190500000001;  I=0.3864255740794268; J=0.3864255740794268;
            M=8.022733000000001;  N=24.068199;           O=0.6237293899283516;  P=0.6237293899283515;
            Q=163.2944534287884;  R=1469.650080859095;   S=0.2410249875732739;  T=0.2410249875732739 }
    9  = @{ G=20.35396833333334;  H=61.06190500000001;  I=0.3864255740794268; J=0.3864255740794268;
            M=0.2404413333333334; N=0.7213240000000001;  O=0.0186931717849216;  P=0.0186931717849216;
            Q=4.893935284691113;  R=44.04541756222001;   S=0.007223519638353673; T=0.007223519638353671 }
    10 = @{ G=24.47399366666667;  H=73.421981;           I=0.4646453653546145; J=0.4646453653546145;
            M=1.115279333333333;  N=3.345838;            O=0.08670767158519405; P=0.08670767158519403;
            Q=27.29533934056422;  R=245.658054065078;    S=0.04028831774275041; T=0.04028831774275041 }
    11 = @{ G=24.47399366666667;  H=73.421981;           I=0.4646453653546145; J=0.4646453653546145;
            M=3.484068333333333;  N=10.452205;           O=0.2708697667015328;  P=0.2708697667015328;
            Q=85.26906632423389;  R=767.421596918105;    S=0.1258583817125529;  T=0.1258583817125529 }
    12 = @{ G=24.47399366666667;  H=73.421981;           I=0.4646453653546145; J=0.4646453653546145;
            M=8.022733000000001;  N=24.068199;           O=0.6237293899283516;  P=0.6237293899283515;
            Q=196.3483166313577;  R=1767.134849682219;   S=0.2898129702656697;  T=0.2898129702656697 }
    13 = @{ G=24.47399366666667;  H=73.421981;           I=0.4646453653546145; J=0.4646453653546145;
            M=0.2404413333333334; N=0.7213240000000001;  O=0.0186931717849216;  P=0.0186931717849216;
            Q=5.88455966920489;   R=52.96103702284401;   S=0.008685695633641468; T=0.008685695633641468 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}

# Remove the "Resolving-Mac" sending-cluster rows (previously rows 14-17) entirely.
$ws.Range("A14:T17").EntireRow.Delete()
